$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column B (Campo Magnético (T) -> Campo Magnético (mT))
$ws.Range("B1").Value = "Campo Magnético (mT)"

# New error columns
$ws.Range("C1").Value = "error_corriente (A)"
$ws.Range("D1").Value = "error_campo (mT)"

# Fill error values for data rows 2-59
$ws.Range("C2:C59").Value = 0.01
$ws.Range("D2:D59").Value = 0.1

# Column widths to match bestFit sizing from the target workbook
$ws.Columns.Item(2).ColumnWidth = 19.44140625
$ws.Columns.Item(3).ColumnWidth = 15.6640625
$ws.Columns.Item(4).ColumnWidth = 15.21875
